$d = $word.ActiveDocument

# 1. Merge the "e.g." runs into one run and remove the proofErr wrapping.
$d.Content.Find.Execute(". Test cases generation should be ahead of the schedule by one week (e.g. ", $true, $false, $false, $false, $false, $true, 1, $false, ". Test cases generation should be ahead of the schedule by one week (e.g. ", 2) | Out-Null

# 2. Remove " [Iteration 1 code submission: Recess Week Mon]" after "Week 6"
$d.Content.Find.Execute(" [Iteration 1 code submission: Recess Week Mon]", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 3. Remove " [Iteration 1 report draft submission: Before Week 7 Consultation]" after "Recess Week"
$d.Content.Find.Execute(" [Iteration 1 report draft submission: Before Week 7 Consultation]", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null
